$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
# Prefixing D-column numeric-looking values with an apostrophe forces Excel
# to keep them as literal text (matching the original "inlineStr" cells)
# instead of auto-converting them into numbers; resetting the Style back to
# "Normal" afterwards avoids leaving a stray quote-prefix number format behind.
$ws.Range("D2").Value = "'96.900.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "'3.674.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'240.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  +11.14%  "
$ws.Range("D7").Value = "'657.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.428"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'3.675.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "'45.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("D15").Value = "'4.358.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("E16").Value = "  +5.02%  "
$ws.Range("D17").Value = "'96.616.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "'8.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").Value = "'3.680.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "'18.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.86%  "
$ws.Range("D21").Value = "'12.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'0.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'533.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "'102.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'13.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("D30").Value = "'12.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.49%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  +17.07%  "
$ws.Range("D34").Value = "'0.187"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "'663.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.06%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'32.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  +4.93%  "
$ws.Range("D39").Value = "'8.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").Value = "'6.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.69%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").Value = "'38.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.09%  "
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("D47").Value = "'0.434"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.96%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'8.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.57%  "

# --- Rows 48 and 49 swap order (Stacks now ranks above MantraDAO) and refresh values ---
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "'3.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.15%  "
